$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update C10 (the "From" value for rule R20) from 18 to 1.
$ws.Range("C10").Value = 1
